$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.709.15"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.022.40"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.60"
$ws.Range("E5").Value = "  -9.89%  "
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.99"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.93"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "2.314.57"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.20"
$ws.Range("E15").Value = "  -6.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.760"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "2.019.66"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "36.777.83"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.85"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.35"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.32"
$ws.Range("E23").Value = "  -5.45%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  -8.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.38"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.91"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0604"
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.75"
$ws.Range("E40").Value = "  +4.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").Value = "1.464.42"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0204"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.35"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.11"
$ws.Range("E46").Value = "  -7.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.32"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.83"
$ws.Range("E49").Value = "  +26.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -2.44%  "
